# Combined PRD for ATR 72-600 and A220-100 showing flight data for US market
#
# 1) Rename existing "ATR 72-600" sheet to "ATR 72-600_pax" (keeps its
#    original data - it was the passenger/pax variant all along).
# 2) Insert a brand-new "ATR 72-600" sheet (the new aircraft-range PRD)
#    right after it.
# 3) Insert a brand-new "A220-100" sheet right after "A220-300".
# 4) Populate the two new sheets with their flight-data tables.
# 5) Add a threaded review comment on A1 of every data sheet citing the
#    source of the digitized data.
# 6) Restore cell selections / active sheet to match the saved view state.

$wb = $excel.ActiveWorkbook

# --- 1) Rename first sheet -------------------------------------------------
$sheetAtrPax = $wb.Worksheets.Item(1)
$sheetAtrPax.Name = "ATR 72-600_pax"

# --- 2) New "A220-100" sheet right after "A220-300" -------------------------
# (created first so the internal sheetId sequence matches the authored file:
#  A220-100 picks up sheetId 4, ATR 72-600 picks up sheetId 5)
$sheetA220300 = $wb.Worksheets.Item("A220-300")
$sheetA220100 = $wb.Worksheets.Add($null, $sheetA220300)
$sheetA220100.Name = "A220-100"

$sheetA220100.Range("A1").Value = "x"
$sheetA220100.Range("B1").Value = "y"

$sheetA220100.Range("A2").Value = 0
$sheetA220100.Range("B2").Value = 115.992957746478

$sheetA220100.Range("A3").Value = 956.15514333895396
$sheetA220100.Range("B3").Value = 115.948356807511

$sheetA220100.Range("A4").Value = 2098.2293423271499
$sheetA220100.Range("B4").Value = 115.45774647887301

$sheetA220100.Range("A5").Value = 3430.0168634063998
$sheetA220100.Range("B5").Value = 104.530516431924

$sheetA220100.Range("A6").Value = 3771.50084317032
$sheetA220100.Range("B6").Value = 101.85446009389599

$sheetA220100.Range("A7").Value = 4435.4974704890301
$sheetA220100.Range("B7").Value = 80.044600938967093

# --- 3) New "ATR 72-600" sheet right after "ATR 72-600_pax" ----------------
$sheetAtr = $wb.Worksheets.Add($null, $sheetAtrPax)
$sheetAtr.Name = "ATR 72-600"

$sheetAtr.Range("A1").Value = "x"
$sheetAtr.Range("B1").Value = "y"

$sheetAtr.Range("A2").Value = 2.8571428571427702
$sheetAtr.Range("B2").Value = 7323.1132075471696
$sheetAtr.Range("D2").Value = 2.8571428571427702
$sheetAtr.Range("E2").Value = 7323.1132075471696

$sheetAtr.Range("A3").Value = 408.57142857142799
$sheetAtr.Range("B3").Value = 7334.9056603773597
$sheetAtr.Range("D3").Value = 408.57142857142799
$sheetAtr.Range("E3").Value = 7334.9056603773597

$sheetAtr.Range("A4").Value = 825.71428571428498
$sheetAtr.Range("B4").Value = 6344.3396226415098
$sheetAtr.Range("D4").Value = 739.99999999999898
$sheetAtr.Range("E4").Value = 6544.8113207547103

$sheetAtr.Range("A5").Value = 1639.99999999999
$sheetAtr.Range("B5").Value = 4375
$sheetAtr.Range("D5").Value = 1639.99999999999
$sheetAtr.Range("E5").Value = 4375

$sheetAtr.Range("A6").Value = 1751.42857142857
$sheetAtr.Range("B6").Value = 23.5849056603783
$sheetAtr.Range("D6").Value = 1751.42857142857
$sheetAtr.Range("E6").Value = 23.5849056603783

# --- 4) Threaded source comments on A1 of each data sheet -------------------
$sheetAtrPax.Range("A1").AddCommentThreaded("Data from figure 9 in high_whitepaper_2023_perFCretroac_mukhopadhaya, see `"C:\Users\nmb48\Documents\GitHub\Flydrogen\misc\fyr\FYR_DigitGraph_Data\fig9_mukhopadhaya_LH2_yellow.csv`"") | Out-Null

$sheetAtr.Range("A1").AddCommentThreaded("Data from figure 79 in `"C:\Users\nmb48\OneDrive - University of Cambridge\Desktop\PhD\Literature\2nd Year\Reports\Restricted FlyZero Reports\FZO-AIN-REP-0008 - Regional Aircraft Concept Report.pdf`"") | Out-Null

$sheetA220300.Range("A1").AddCommentThreaded("Data from page 223 in A220-ACP-Issue011-00-18Sep2025") | Out-Null

$sheetA220100.Range("A1").AddCommentThreaded("Data from page 209 in A220-ACP-Issue011-00-18Sep2025") | Out-Null

# --- 5) Selections matching the saved view state -----------------------------
$sheetAtrPax.Activate()
$sheetAtrPax.Range("K9").Select()

$sheetA220100.Activate()
$sheetA220100.Range("A3").Select()

$sheetA220300.Activate()

$sheetAtr.Activate()
$sheetAtr.Range("K20").Select()
